# Applies the commit's change: several lines get their text split across
# multiple runs separated by <w:proofErr w:type="gramStart"/>.../<w:proofErr
# w:type="gramEnd"/> pairs (as Word's grammar checker would insert while the
# text was retyped), and three additional blank paragraphs are inserted in
# the gap between the two YAML documents.

$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# Single pass over the (unchanging, for this step) paragraph collection to
# find every index we need to touch -- avoids re-scanning the whole
# document once per replacement.
$idxSharedWorkspaceName = New-Object System.Collections.ArrayList
$idxWorkspaceShared = -1
$idxInstallReq = -1
$idxSteps = -1
$idxWorkingDir = -1
$idxBash = -1
$idxPipInstall = -1

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)
    switch ($t) {
        "    - name: shared-workspace"                          { [void]$idxSharedWorkspaceName.Add($i) }
        "          workspace: shared-workspace"                 { $idxWorkspaceShared = $i }
        "    - name: install-requirements"                      { $idxInstallReq = $i }
        "        steps:"                                        { $idxSteps = $i }
        "            workingDir: /workspace/shared-workspace"   { $idxWorkingDir = $i }
        "              #!/usr/bin/env bash"                     { $idxBash = $i }
        "              pip install -r requirements.txt"         { $idxPipInstall = $i }
    }
}

# ---------------------------------------------------------------------
# 1) "    - name: shared-workspace"  (appears twice, identical paragraphs)
#    -> "    - name: " + gramStart + "shared-workspace" + gramEnd
# ---------------------------------------------------------------------
$xmlSharedWorkspaceName = "<w:p $wns><w:r><w:t xml:space='preserve'>    - name: </w:t></w:r><w:proofErr w:type='gramStart'/><w:r><w:t>shared-workspace</w:t></w:r><w:proofErr w:type='gramEnd'/></w:p>"
foreach ($i in $idxSharedWorkspaceName) {
    [void]$d.Paragraphs.Item($i).Range.InsertXML($xmlSharedWorkspaceName)
}

# ---------------------------------------------------------------------
# 2) "          workspace: shared-workspace"
#    -> "          workspace: " + gramStart + "shared-workspace" + gramEnd
# ---------------------------------------------------------------------
if ($idxWorkspaceShared -ne -1) {
    $xml = "<w:p $wns><w:r><w:t xml:space='preserve'>          workspace: </w:t></w:r><w:proofErr w:type='gramStart'/><w:r><w:t>shared-workspace</w:t></w:r><w:proofErr w:type='gramEnd'/></w:p>"
    [void]$d.Paragraphs.Item($idxWorkspaceShared).Range.InsertXML($xml)
}

# ---------------------------------------------------------------------
# 3) "    - name: install-requirements"
#    -> "    - name: " + gramStart + "install-requirements" + gramEnd
# ---------------------------------------------------------------------
if ($idxInstallReq -ne -1) {
    $xml = "<w:p $wns><w:r><w:t xml:space='preserve'>    - name: </w:t></w:r><w:proofErr w:type='gramStart'/><w:r><w:t>install-requirements</w:t></w:r><w:proofErr w:type='gramEnd'/></w:p>"
    [void]$d.Paragraphs.Item($idxInstallReq).Range.InsertXML($xml)
}

# ---------------------------------------------------------------------
# 4) "        steps:"
#    -> "        " + gramStart + "steps" + gramEnd + ":"
# ---------------------------------------------------------------------
if ($idxSteps -ne -1) {
    $xml = "<w:p $wns><w:r><w:t xml:space='preserve'>        </w:t></w:r><w:proofErr w:type='gramStart'/><w:r><w:t>steps</w:t></w:r><w:proofErr w:type='gramEnd'/><w:r><w:t>:</w:t></w:r></w:p>"
    [void]$d.Paragraphs.Item($idxSteps).Range.InsertXML($xml)
}

# ---------------------------------------------------------------------
# 5) "            workingDir: /workspace/shared-workspace"
#    (proofErr spellStart/workingDir/spellEnd stays; only the trailing run
#     ": /workspace/shared-workspace" is split)
#    -> ": /workspace/" + gramStart + "shared-workspace" + gramEnd
# ---------------------------------------------------------------------
if ($idxWorkingDir -ne -1) {
    $xml = "<w:p $wns><w:r><w:lastRenderedPageBreak/><w:t xml:space='preserve'>            </w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:t>workingDir</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t>: /workspace/</w:t></w:r><w:proofErr w:type='gramStart'/><w:r><w:t>shared-workspace</w:t></w:r><w:proofErr w:type='gramEnd'/></w:p>"
    [void]$d.Paragraphs.Item($idxWorkingDir).Range.InsertXML($xml)
}

# ---------------------------------------------------------------------
# 6) "              #!/usr/bin/env bash"
#    -> "              " + gramStart + "#!/" + gramEnd + "usr/bin/env bash"
# ---------------------------------------------------------------------
if ($idxBash -ne -1) {
    $xml = "<w:p $wns><w:r><w:t xml:space='preserve'>              </w:t></w:r><w:proofErr w:type='gramStart'/><w:r><w:t>#!/</w:t></w:r><w:proofErr w:type='gramEnd'/><w:r><w:t>usr/bin/env bash</w:t></w:r></w:p>"
    [void]$d.Paragraphs.Item($idxBash).Range.InsertXML($xml)
}

# ---------------------------------------------------------------------
# 7) Insert three extra blank paragraphs right after the
#    "pip install -r requirements.txt" line, before the run of blank
#    paragraphs that separates the two YAML documents. Inserting
#    "<w:p/><w:p/><w:p/>" XML at the collapsed end-of-paragraph point
#    yields clean empty paragraphs (no stray empty run), matching the
#    diff's plain "<w:p/>" additions.
# ---------------------------------------------------------------------
if ($idxPipInstall -ne -1) {
    $endOfPara = $d.Paragraphs.Item($idxPipInstall).Range.End
    $insertPoint = $d.Range($endOfPara, $endOfPara)
    $blankParasXml = "<w:p $wns/><w:p $wns/><w:p $wns/>"
    [void]$insertPoint.InsertXML($blankParasXml)
}
